# Update NATMI LR-pairs output (Efna3-Epha4) with newly recomputed TPM
# values. The sending cluster changes from MuSCs to ECs, the
# "Resolving-Mac" target row is dropped (data no longer produced for that
# cluster), and every numeric column for the remaining four rows is
# refreshed with the new TPM-derived figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-5 (row 6 / Resolving-Mac is removed entirely).
# Columns: A..T
#  A Sending cluster   B Ligand symbol   C Receptor symbol  D Target cluster
#  E..T numeric metrics
$rows = @(
    @{ Row = 2;  A = "ECs"; D = "ECs";               E=1; F=0.3333333333333333; G=0.02551366666666667;  H=0.076541;  I=1; J=1; K=3; L=1;                  M=4.886733666666667;  N=14.660201; O=0.4081653954827624;   P=0.4171266852711343;  Q=0.1246784938601111;   R=1.122106444741;     S=0.4081653954827624;   T=0.4171266852711343 }
    @{ Row = 3;  A = "ECs"; D = "FAPs";               E=1; F=0.3333333333333333; G=0.02551366666666667;  H=0.076541;  I=1; J=1; K=3; L=1;                  M=6.292848333333333;  N=18.878545; O=0.5256114009667484;   P=0.5371512231375235;  Q=0.1605536347605555;   R=1.444982712845;     S=0.5256114009667484;   T=0.5371512231375235 }
    @{ Row = 4;  A = "ECs"; D = "Inflammatory-Mac";   E=1; F=0.3333333333333333; G=0.02551366666666667;  H=0.076541;  I=1; J=1; K=1; L=0.3333333333333333; M=0.021228;           N=0.063684;  O=0.001773072896198643; P=0.001812000792131494; Q=0.000541604116;       R=0.004874437044;     S=0.001773072896198643; T=0.001812000792131494 }
    @{ Row = 5;  A = "ECs"; D = "MuSCs";              E=1; F=0.3333333333333333; G=0.02551366666666667;  H=0.076541;  I=1; J=1; K=2; L=1;                  M=0.771625;           N=1.54325;   O=0.06445013065429045;  P=0.04391009079921059;  Q=0.01968698304166667;  R=0.11812189825;      S=0.06445013065429045;  T=0.04391009079921059 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = "Efna3"
    $ws.Cells.Item($row, 3).Value = "Epha4"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}

# The old row 6 (MuSCs -> Resolving-Mac) is no longer part of the output;
# delete it and shift everything below it up.
$ws.Rows.Item(6).Delete()
